# Fruta / hortaliza, semanal
#
# Insert one new weekly price-report pair (Primera/Segunda) at the top of
# the "Acelga" data block (rows 294-295), pushing the existing rows
# 294-315 down by two rows (to 296-317).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 294 - this shifts the old rows 294:315
# down to 296:317 and keeps per-cell formatting (e.g. the date style on
# column D) consistent with the surrounding rows.
$ws.Range("A294:R295").EntireRow.Insert()

# New row 294 - "Primera" quality entry for the new week.
$ws.Range("A294").Value = 8
$ws.Range("B294").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C294").Value = 'Coquimbo'
$ws.Range("D294").Value = 44585
$ws.Range("E294").Value = 4
$ws.Range("F294").Value = 100112009
$ws.Range("G294").Value = 'Acelga'
$ws.Range("H294").Value = 'Sin especificar'
$ws.Range("I294").Value = 'Primera'
$ws.Range("J294").Value = 2400
$ws.Range("K294").Value = 450
$ws.Range("L294").Value = 500
$ws.Range("M294").Value = 475
$ws.Range("N294").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("O294").Value = 'Provincia del Elquí'
$ws.Range("P294").Value = 238
$ws.Range("Q294").Value = 2
$ws.Range("R294").Value = 'Hortaliza'

# New row 295 - "Segunda" quality entry for the new week.
$ws.Range("A295").Value = 8
$ws.Range("B295").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C295").Value = 'Coquimbo'
$ws.Range("D295").Value = 44585
$ws.Range("E295").Value = 4
$ws.Range("F295").Value = 100112009
$ws.Range("G295").Value = 'Acelga'
$ws.Range("H295").Value = 'Sin especificar'
$ws.Range("I295").Value = 'Segunda'
$ws.Range("J295").Value = 1500
$ws.Range("K295").Value = 350
$ws.Range("L295").Value = 400
$ws.Range("M295").Value = 375
$ws.Range("N295").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("O295").Value = 'Provincia del Elquí'
$ws.Range("P295").Value = 188
$ws.Range("Q295").Value = 2
$ws.Range("R295").Value = 'Hortaliza'
